$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: "Authorship Resource" header + value for each data row
$ws.Range("K1").Value = "Authorship Resource"
$ws.Range("K2:K6").Value = "Noémi Villars-Amberg, Daniela Subotic"

# Match the cursor / selection state recorded in the saved file
$ws.Range("K17").Select()
